$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the XPath header strings in row 1 (C1:G1): double quotes -> single quotes
$ws.Range("C1").Value = "//*[@id='insurance-form']/div/section[1]"
$ws.Range("D1").Value = "//*[@id='insurance-form']/div/section[2]"
$ws.Range("E1").Value = "//*[@id='insurance-form']/div/section[3]"
$ws.Range("F1").Value = "//*[@id='insurance-form']/div/section[4]"
$ws.Range("G1").Value = "//*[@id='insurance-form']/div/section[5]"

# Move the selection/active cell to F29:F30
$ws.Range("F29:F30").Select() | Out-Null
